# Insert a new weekly record at row 130, shifting all existing rows
# from 130 downward to 131 onward (matches the commit: "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 130; existing row 130 (and everything
# below it) moves down to row 131, etc.
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row 130 with the new weekly observation.
$ws.Range("A130").Value = 8
$ws.Range("B130").Value = "Terminal La Palmera de La Serena"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 45167
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 100112040
$ws.Range("G130").Value = "Cilantro"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 2000
$ws.Range("K130").Value = 2000
$ws.Range("L130").Value = 2500
$ws.Range("M130").Value = 2250
$ws.Range("N130").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O130").Value = "Provincia del Elquí"
$ws.Range("P130").Value = 1500
$ws.Range("Q130").Value = 1.5
$ws.Range("R130").Value = "Hortaliza"
